# Issue #20 images to playlist return nav
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Row 22: new BUG issue #21 - adding all images also adds the directory itself
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 1
$ws.Range("D22").Value = "BUG"
$ws.Range("H22").Value = "When adding all images it also adds the directory"
$ws.Range("E22").Value = "don’t add directory to playlist"
$ws.Rows.Item(22).RowHeight = 29

# Row 23: new BUG issue #22 - playlist save not working
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 1
$ws.Range("D23").Value = "BUG"
$ws.Range("H23").Value = "Playlist save not working"

# Leave the cursor where the editor last left it (scrolled down a little)
$ws.Range("E29").Select()
